$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row with the VIVA problem placeholder text
$ws.Range("A7").Value = "????"

# Update the selected cell to mirror the author's saved view state
$ws.Range("A17").Select()
